# Applies the "Added initial simulation and verification slides" edit:
#  - inserts 3 new slides (Major Design Decisions, Dynamic Model, Kinematic Modelling)
#  - renames/fills in content on existing placeholder slides
#  - updates the "Background" slide into "Project Overview" with bullet content

$p = $ppt.ActivePresentation

# NOTE: this PowerShell host does not honour named (-Param value) arguments
# on user-defined functions, so Set-Body is called positionally everywhere.
function Set-Body {
    param($Shape, $Lines, $Levels)
    $tr = $Shape.TextFrame.TextRange
    $tr.Text = [string]::Join("`r", $Lines)
    for ($i = 0; $i -lt $Levels.Count; $i++) {
        if ($Levels[$i] -gt 1) {
            $tr.Paragraphs($i + 1).IndentLevel = $Levels[$i]
        }
    }
}

# ---------------------------------------------------------------------------
# Slide 2 (existing): "Background" -> "Project Overview" with new bullet body
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Project Overview"

$lines2 = @(
    "Develop initial robotics platform",
    "Public Outreach",
    "Increase interest in STEM fields",
    "Engineering deficit",
    "Educational robotics platform",
    "Fluid Power",
    "Control Systems"
)
$levels2 = @(1, 1, 2, 2, 1, 2, 2)
Set-Body $s2.Shapes.Item(2) $lines2 $levels2

# ---------------------------------------------------------------------------
# New slide at position 3: "Major Design Decisions"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Major Design Decisions"

$lines3 = @(
    "Quadruped Locomotion",
    "Agile motion",
    "Stable",
    "Rugged Terrain",
    "Pneumatic Power",
    "High energy density",
    "Clean and lightweight"
)
$levels3 = @(1, 2, 2, 2, 1, 2, 2)
Set-Body $s3.Shapes.Item(2) $lines3 $levels3

# ---------------------------------------------------------------------------
# Existing slides 4-6 keep their (empty) bodies, just get new titles so the
# deck reads: Mechanical Design, Pneumatic Design, Simulation and Verification
# ---------------------------------------------------------------------------
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "Mechanical Design"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "Pneumatic Design"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "Simulation and Verification"

# ---------------------------------------------------------------------------
# New slide at position 7: "Dynamic Model"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Dynamic Model"

$lines7 = @(
    "Lagrangian Mechanics",
    "Non-inertial reference frames",
    "Arbitrary grounded reference frame",
    "Determines mechanical properties of components",
    "27 Equations"
)
$levels7 = @(1, 1, 1, 1, 1)
Set-Body $s7.Shapes.Item(2) $lines7 $levels7

# ---------------------------------------------------------------------------
# New slide at position 8: "Kinematic Modelling"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Kinematic Modelling"

$lines8 = @(
    "Homogenous Transformations",
    "Determines cylinder stroke"
)
$levels8 = @(1, 1)
Set-Body $s8.Shapes.Item(2) $lines8 $levels8

# ---------------------------------------------------------------------------
# Existing slides 9-10 (previously "Control architecture and electronics"
# and "Questions?") keep their text/position at the end of the deck already,
# nothing further to change there.
# ---------------------------------------------------------------------------

Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $title = $s.Shapes.Item(1).TextFrame.TextRange.Text
    Write-Output "Slide $i (id=$($s.SlideID)): $title"
}
